$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47
$ws.Range("A47").Value = 111736272
$ws.Range("B47").Value = 89401
$ws.Range("E47").Value = 1108
$ws.Range("F47").Value = 'Harticka'
$ws.Range("G47").Value = 'Pelloporus leporinus'
$ws.Range("H47").Value = '(Fr.) Krieglst.'
$ws.Range("M47").Value = $null
$ws.Range("Q47").Value = 616308.8236423519
$ws.Range("R47").Value = 7268903.133137755

# Row 48
$ws.Range("A48").Value = 111736525
$ws.Range("B48").Value = 56398
$ws.Range("E48").Value = 100109
$ws.Range("F48").Value = 'Tretåig hackspett'
$ws.Range("G48").Value = 'Picoides tridactylus'
$ws.Range("H48").Value = '(Linnaeus, 1758)'
$ws.Range("M48").Value = 'färska spår'
$ws.Range("Q48").Value = 616358.6131022752
$ws.Range("R48").Value = 7268822.486957001
$ws.Range("S48").Value = 25

# Row 49
$ws.Range("A49").Value = 111736370
$ws.Range("B49").Value = 56398
$ws.Range("E49").Value = 100109
$ws.Range("F49").Value = 'Tretåig hackspett'
$ws.Range("G49").Value = 'Picoides tridactylus'
$ws.Range("H49").Value = '(Linnaeus, 1758)'
$ws.Range("M49").Value = 'färska spår'
$ws.Range("Q49").Value = 616327.1020967637
$ws.Range("R49").Value = 7268872.304318298

# Row 50
$ws.Range("A50").Value = 111736257
$ws.Range("B50").Value = 77515
$ws.Range("E50").Value = 6425
$ws.Range("F50").Value = 'Garnlav'
$ws.Range("G50").Value = 'Alectoria sarmentosa'
$ws.Range("H50").Value = '(Ach.) Ach.'

# Row 51
$ws.Range("A51").Value = 111736405
$ws.Range("B51").Value = 77515
$ws.Range("E51").Value = 6425
$ws.Range("F51").Value = 'Garnlav'
$ws.Range("G51").Value = 'Alectoria sarmentosa'
$ws.Range("H51").Value = '(Ach.) Ach.'
$ws.Range("M51").Value = $null
$ws.Range("Q51").Value = 616333.1117616051
$ws.Range("R51").Value = 7268857.179896916
$ws.Range("S51").Value = 10
$ws.Range("AC51").Value = $null

# Row 52
$ws.Range("A52").Value = 111736402
$ws.Range("B52").Value = 89423
$ws.Range("E52").Value = 5432
$ws.Range("F52").Value = 'Granticka'
$ws.Range("G52").Value = 'Porodaedalea chrysoloma'
$ws.Range("H52").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q52").Value = 616333.1117616051
$ws.Range("R52").Value = 7268857.179896916

# Row 53
$ws.Range("A53").Value = 111778126
$ws.Range("B53").Value = 89405
$ws.Range("E53").Value = 1202
$ws.Range("F53").Value = 'Ullticka'
$ws.Range("G53").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H53").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("M53").Value = $null
$ws.Range("Q53").Value = 616202.3044715263
$ws.Range("R53").Value = 7268603.611313918
$ws.Range("AC53").Value = $null

# Row 54
$ws.Range("A54").Value = 111778005
$ws.Range("B54").Value = 89369
$ws.Range("D54").Value = 'LC'
$ws.Range("E54").Value = 5447
$ws.Range("F54").Value = 'Vedticka'
$ws.Range("G54").Value = 'Fuscoporia viticola'
$ws.Range("H54").Value = '(Schwein.) Murrill'
$ws.Range("Q54").Value = 616499.3130462242
$ws.Range("R54").Value = 7268610.508796399

# Row 55
$ws.Range("A55").Value = 111777494
$ws.Range("B55").Value = 90854
$ws.Range("E55").Value = 2079
$ws.Range("F55").Value = 'Nordtagging'
$ws.Range("G55").Value = 'Odonticium romellii'
$ws.Range("H55").Value = '(S.Lundell) Parmasto'
$ws.Range("M55").Value = $null
$ws.Range("AC55").Value = $null

# Row 56
$ws.Range("A56").Value = 111777467
$ws.Range("B56").Value = 89423
$ws.Range("E56").Value = 5432
$ws.Range("F56").Value = 'Granticka'
$ws.Range("G56").Value = 'Porodaedalea chrysoloma'
$ws.Range("H56").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q56").Value = 616413.4864248879
$ws.Range("R56").Value = 7268760.315060399

# Row 57
$ws.Range("A57").Value = 111777491
$ws.Range("B57").Value = 56398
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = 'Tretåig hackspett'
$ws.Range("G57").Value = 'Picoides tridactylus'
$ws.Range("H57").Value = '(Linnaeus, 1758)'
$ws.Range("M57").Value = 'färska spår'
$ws.Range("AC57").Value = 'Skalade granstammar'

# Row 58
$ws.Range("A58").Value = 111777447
$ws.Range("B58").Value = 89405
$ws.Range("E58").Value = 1202
$ws.Range("F58").Value = 'Ullticka'
$ws.Range("G58").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H58").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("M58").Value = $null
$ws.Range("Q58").Value = 616379.7321599644
$ws.Range("R58").Value = 7268803.814155157
$ws.Range("AC58").Value = $null

# Row 59
$ws.Range("A59").Value = 111778163
$ws.Range("B59").Value = 56398
$ws.Range("D59").Value = 'NT'
$ws.Range("E59").Value = 100109
$ws.Range("F59").Value = 'Tretåig hackspett'
$ws.Range("G59").Value = 'Picoides tridactylus'
$ws.Range("H59").Value = '(Linnaeus, 1758)'
$ws.Range("M59").Value = 'färska spår'
$ws.Range("Q59").Value = 616207.2556492372
$ws.Range("R59").Value = 7268635.7870906
$ws.Range("S59").Value = 10
$ws.Range("AC59").Value = 'Skalad gran'

# Row 60
$ws.Range("A60").Value = 111777411
$ws.Range("B60").Value = 56398
$ws.Range("E60").Value = 100109
$ws.Range("F60").Value = 'Tretåig hackspett'
$ws.Range("G60").Value = 'Picoides tridactylus'
$ws.Range("H60").Value = '(Linnaeus, 1758)'
$ws.Range("M60").Value = 'färska spår'
$ws.Range("Q60").Value = 616367.7277224116
$ws.Range("R60").Value = 7268802.503264537
$ws.Range("AC60").Value = 'Skalade granstammar'

# Row 61
$ws.Range("A61").Value = 111778248
$ws.Range("B61").Value = 56398
$ws.Range("D61").Value = 'NT'
$ws.Range("E61").Value = 100109
$ws.Range("F61").Value = 'Tretåig hackspett'
$ws.Range("G61").Value = 'Picoides tridactylus'
$ws.Range("H61").Value = '(Linnaeus, 1758)'
$ws.Range("M61").Value = 'färsk spillning'
$ws.Range("Q61").Value = 616162.9874832245
$ws.Range("R61").Value = 7268630.281087617
$ws.Range("AC61").Value = 'Skalad gran'

# Row 62
$ws.Range("A62").Value = 111777380
$ws.Range("B62").Value = 56398
$ws.Range("E62").Value = 100109
$ws.Range("F62").Value = 'Tretåig hackspett'
$ws.Range("G62").Value = 'Picoides tridactylus'
$ws.Range("H62").Value = '(Linnaeus, 1758)'
$ws.Range("M62").Value = 'färska spår'
$ws.Range("Q62").Value = 616414.0528149965
$ws.Range("R62").Value = 7268860.418718725
$ws.Range("AC62").Value = 'Skalade stammar'

# Row 63
$ws.Range("A63").Value = 111777940
$ws.Range("B63").Value = 90678
$ws.Range("D63").Value = 'LC'
$ws.Range("E63").Value = 4366
$ws.Range("F63").Value = 'Skarp dropptaggsvamp'
$ws.Range("G63").Value = 'Hydnellum peckii'
$ws.Range("H63").Value = 'Banker'
$ws.Range("M63").Value = $null
$ws.Range("Q63").Value = 616438.7745429112
$ws.Range("R63").Value = 7268803.685732875
$ws.Range("S63").Value = 25
$ws.Range("AC63").Value = $null

# Row 64
$ws.Range("A64").Value = 111777499
$ws.Range("B64").Value = 78107
$ws.Range("E64").Value = 6453
$ws.Range("F64").Value = 'Vedskivlav'
$ws.Range("G64").Value = 'Hertelidea botryosa'
$ws.Range("H64").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("M64").Value = $null
$ws.Range("Q64").Value = 616426.5202303537
$ws.Range("R64").Value = 7268746.301918368
$ws.Range("AC64").Value = $null

# Row 65
$ws.Range("A65").Value = 111777331
$ws.Range("B65").Value = 89423
$ws.Range("E65").Value = 5432
$ws.Range("F65").Value = 'Granticka'
$ws.Range("G65").Value = 'Porodaedalea chrysoloma'
$ws.Range("H65").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q65").Value = 616362.7639770868
$ws.Range("R65").Value = 7268822.653031595
